# Update Avp-Avpr2.xlsx data:
#  - Row 2 (MuSCs -> Avp -> Avpr2 -> Inflammatory-Mac): refresh TPM-derived metrics
#  - Row 3 (MuSCs -> Avp -> Avpr2 -> Neutrophils): was Resolving-Mac row, now becomes the
#    MuSCs->Neutrophils row with refreshed metrics
#  - Row 4 (MuSCs -> Avp -> Avpr2 -> Resolving-Mac): was first Neutrophils row, now becomes
#    the MuSCs->Resolving-Mac row with refreshed metrics
#  - Old row 5 (Neutrophils -> Avp -> Avpr2 -> Resolving-Mac) is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Avp"
$ws.Range("C2").Value = "Avpr2"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.2600935
$ws.Range("H2").Value = 0.520187
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1811646666666666
$ws.Range("N2").Value = 0.5434939999999999
$ws.Range("O2").Value = 0.3599820106359796
$ws.Range("P2").Value = 0.3599820106359796
$ws.Range("Q2").Value = 0.04711975222966665
$ws.Range("R2").Value = 0.2827185133779999
$ws.Range("S2").Value = 0.3599820106359796
$ws.Range("T2").Value = 0.3599820106359796

# --- Row 3 ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Avp"
$ws.Range("C3").Value = "Avpr2"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.2600935
$ws.Range("H3").Value = 0.520187
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.234272
$ws.Range("N3").Value = 0.702816
$ws.Range("O3").Value = 0.4655085737600355
$ws.Range("P3").Value = 0.4655085737600355
$ws.Range("Q3").Value = 0.060932624432
$ws.Range("R3").Value = 0.3655957465919999
$ws.Range("S3").Value = 0.4655085737600355
$ws.Range("T3").Value = 0.4655085737600355

# --- Row 4 ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Avp"
$ws.Range("C4").Value = "Avpr2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.2600935
$ws.Range("H4").Value = 0.520187
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08782366666666667
$ws.Range("N4").Value = 0.263471
$ws.Range("O4").Value = 0.174509415603985
$ws.Range("P4").Value = 0.174509415603985
$ws.Range("Q4").Value = 0.02284236484616667
$ws.Range("R4").Value = 0.137054189077
$ws.Range("S4").Value = 0.174509415603985
$ws.Range("T4").Value = 0.174509415603985

# --- Remove old row 5 (Neutrophils -> Avp -> Avpr2 -> Resolving-Mac) ---
$ws.Range("A5:T5").EntireRow.Delete()
